$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data scraped on Sat Dec 23 18:41:43 UTC 2023

$ws.Range('D2').Value = '43.894.84'
$ws.Range('E2').Value = '  -0.19%  '
$ws.Range('D3').Value = '2.299.36'
$ws.Range('E3').Value = '  -1.43%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = "'99.92"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +2.39%  '
$ws.Range('D6').Value = "'270.66"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.43%  '
$ws.Range('D7').Value = "'0.626"
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.27%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').Value = "'0.608"
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -3.29%  '
$ws.Range('D10').Value = "'44.87"
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -2.50%  '
$ws.Range('D11').Value = "'0.0930"
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -2.20%  '
$ws.Range('E12').Value = '  -2.95%  '
$ws.Range('E13').Value = '  +1.59%  '
$ws.Range('D14').Value = "'15.85"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +1.19%  '
$ws.Range('D15').Value = '2.642.75'
$ws.Range('E15').Value = '  -1.66%  '
$ws.Range('D16').Value = "'0.857"
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -1.66%  '
$ws.Range('D17').Value = '2.293.39'
$ws.Range('E17').Value = '  -1.89%  '
$ws.Range('D18').Value = '43.830.17'
$ws.Range('E18').Value = '  -0.16%  '
$ws.Range('E19').Value = '  +1.20%  '
$ws.Range('D20').Value = "'6.24"
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -3.14%  '
$ws.Range('D21').Value = "'72.35"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.84%  '
$ws.Range('E22').Value = '  +7.45%  '
$ws.Range('D23').Value = "'233.35"
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -2.86%  '
$ws.Range('D24').Value = "'2.88"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +13.44%  '
$ws.Range('D25').Value = "'9.14"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -3.67%  '
$ws.Range('E26').Value = '  +0.09%  '
$ws.Range('D27').Value = "'11.22"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -2.14%  '
$ws.Range('E28').Value = '  -0.89%  '
$ws.Range('D29').Value = "'2.29"
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.11%  '
$ws.Range('D30').Value = "'38.39"
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +0.16%  '
$ws.Range('E31').Value = '  +1.89%  '
$ws.Range('D32').Value = "'21.85"
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -2.88%  '
$ws.Range('D33').Value = "'0.0895"
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -1.27%  '
$ws.Range('D34').Value = "'5.45"
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -1.07%  '
$ws.Range('E35').Value = '  +0.49%  '
$ws.Range('D36').Value = "'4.74"
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +6.67%  '
$ws.Range('E37').Value = '  -1.14%  '
$ws.Range('D38').Value = "'0.0353"
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -2.67%  '
$ws.Range('E39').Value = '  +3.81%  '
$ws.Range('E40').Value = '  -0.96%  '
$ws.Range('D41').Value = "'2.34"
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -1.50%  '
$ws.Range('E42').Value = '  -0.01%  '
$ws.Range('D43').Value = "'12.23"
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.47%  '
$ws.Range('D44').Value = "'64.94"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +3.57%  '
$ws.Range('D45').Value = "'8.85"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -4.33%  '
$ws.Range('D46').Value = "'5.24"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -3.08%  '
$ws.Range('E47').Value = '  -1.60%  '
$ws.Range('E48').Value = '  +0.61%  '
$ws.Range('D49').Value = "'98.62"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -2.00%  '
$ws.Range('B50').Value = 'WOONetwork'
$ws.Range('C50').Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range('D50').Value = "'0.442"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +5.35%  '
$ws.Range('B51').Value = 'Stacks'
$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D51').Value = "'1.53"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +11.23%  '
